# Apply updated cryptocurrency price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "36.531.84"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.959.37"
$ws.Range("E3").Value = "  +1.44%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.18"
$ws.Range("E5").Value = "  +1.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.614"
$ws.Range("E6").Value = "  +1.49%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.69"
$ws.Range("E7").Value = "  +3.71%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.373"
$ws.Range("E9").Value = "  +4.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0789"
$ws.Range("E10").Value = "  -5.16%  "
$ws.Range("E11").Value = "  -1.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.25"
$ws.Range("E12").Value = "  +6.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.839"
$ws.Range("E13").Value = "  +5.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.247.92"
$ws.Range("E14").Value = "  +1.37%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.33"
$ws.Range("E15").Value = "  +2.70%  "
$ws.Range("E16").Value = "  +3.50%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.963.33"
$ws.Range("E17").Value = "  +1.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.494.87"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.69"
$ws.Range("E19").Value = "  +1.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0849"
$ws.Range("E20").Value = "  -0.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "229.74"
$ws.Range("E21").Value = "  +1.39%  "
$ws.Range("E22").Value = "  +2.82%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.45"
$ws.Range("E24").Value = "  +5.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.36"
$ws.Range("E25").Value = "  +4.43%  "
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.140"
$ws.Range("E26").Value = "  +9.16%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.18"
$ws.Range("E27").Value = "  -0.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.72"
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.43"
$ws.Range("E29").Value = "  +1.87%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.119"
$ws.Range("E30").Value = "  +2.19%  "
$ws.Range("E31").Value = "  +8.30%  "
$ws.Range("E32").Value = "  +4.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0613"
$ws.Range("E33").Value = "  -1.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.39"
$ws.Range("E34").Value = "  +6.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.54"
$ws.Range("E35").Value = "  +20.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.30"
$ws.Range("E36").Value = "  +9.18%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("E38").Value = "  -1.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.47"
$ws.Range("E39").Value = "  -9.63%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0978"
$ws.Range("E40").Value = "  +1.17%  "
$ws.Range("E41").Value = "  +1.16%  "
$ws.Range("E42").Value = "  +2.43%  "
$ws.Range("E43").Value = "  +1.30%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.370.82"
$ws.Range("E44").Value = "  +3.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.76"
$ws.Range("E45").Value = "  +2.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.14"
$ws.Range("E46").Value = "  +3.09%  "
$ws.Range("E47").Value = "  +1.66%  "
$ws.Range("E48").Value = "  +1.55%  "
$ws.Range("E49").Value = "  +0.85%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.138.59"
$ws.Range("E50").Value = "  +1.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.00"
$ws.Range("E51").Value = "  +0.76%  "
